# Auto-generated edit script: updates cached numeric values in the
# Leve profit-tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to reflect refreshed currentAveragePrice / LevePrice / LeveProfit figures.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 870.129
$ws.Range("I15").Value = 870.129
$ws.Range("K15").Value = 2610.387
$ws.Range("M15").Value = -2441.387

$ws.Range("H76").Value = 7000.25
$ws.Range("I76").Value = 6000.3335
$ws.Range("K76").Value = 6000.3335
$ws.Range("M76").Value = -5685.3335

$ws.Range("H79").Value = 7000.25
$ws.Range("I79").Value = 6000.3335
$ws.Range("K79").Value = 6000.3335
$ws.Range("M79").Value = -4908.3335

$ws.Range("H88").Value = 3150.75
$ws.Range("J88").Value = 900
$ws.Range("L88").Value = 900
$ws.Range("N88").Value = -1712

$ws.Range("H91").Value = 3150.75
$ws.Range("J91").Value = 900
$ws.Range("L91").Value = 900
$ws.Range("N91").Value = -3708

$ws.Range("H112").Value = 5497369
$ws.Range("J112").Value = 5684658.5
$ws.Range("L112").Value = 17053975.5
$ws.Range("N112").Value = -17056191.5

$ws.Range("H113").Value = 6482.5454
$ws.Range("I113").Value = 5372
$ws.Range("K113").Value = 5372
$ws.Range("M113").Value = -2118

$ws.Range("H116").Value = 7074.1665
$ws.Range("I116").Value = 4499.6665
$ws.Range("J116").Value = 7932.3335
$ws.Range("K116").Value = 4499.6665
$ws.Range("L116").Value = 7932.3335
$ws.Range("M116").Value = -1057.6665
$ws.Range("N116").Value = -14816.3335

$ws.Range("H132").Value = 941.52545
$ws.Range("I132").Value = 864.56366
$ws.Range("J132").Value = 1999.75
$ws.Range("K132").Value = 2593.69098
$ws.Range("L132").Value = 5999.25
$ws.Range("M132").Value = -63.69098000000031
$ws.Range("N132").Value = -11059.25

$ws.Range("H135").Value = 8334786.5
$ws.Range("I135").Value = 1240.05
$ws.Range("J135").Value = 25001880
$ws.Range("K135").Value = 11160.45
$ws.Range("L135").Value = 225016920
$ws.Range("M135").Value = -8625.449999999999
$ws.Range("N135").Value = -225021990

$ws.Range("H137").Value = 46515300
$ws.Range("I137").Value = 34485996
$ws.Range("J137").Value = 71433144
$ws.Range("K137").Value = 103457988
$ws.Range("L137").Value = 214299432
$ws.Range("M137").Value = -103455438
$ws.Range("N137").Value = -214304532

$ws.Range("H138").Value = 5959146
$ws.Range("I138").Value = 2878.7896
$ws.Range("J138").Value = 9017770
$ws.Range("K138").Value = 8636.3688
$ws.Range("L138").Value = 27053310
$ws.Range("M138").Value = -3496.3688
$ws.Range("N138").Value = -27063590

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 615
$ws.Range("I2").Value = 615
$ws.Range("K2").Value = 615
$ws.Range("M2").Value = -502

$ws.Range("H5").Value = 213.11111
$ws.Range("I5").Value = 136
$ws.Range("J5").Value = 367.33334
$ws.Range("K5").Value = 136
$ws.Range("L5").Value = 367.33334
$ws.Range("M5").Value = -24
$ws.Range("N5").Value = -591.33334

$ws.Range("H32").Value = 22232916
$ws.Range("I32").Value = 26322980
$ws.Range("J32").Value = 29711.428
$ws.Range("K32").Value = 26322980
$ws.Range("L32").Value = 29711.428
$ws.Range("M32").Value = -26322693
$ws.Range("N32").Value = -30285.428

$ws.Range("H61").Value = 58828620
$ws.Range("I61").Value = 83336170
$ws.Range("K61").Value = 83336170
$ws.Range("M61").Value = -83335958

$ws.Range("H74").Value = 111237480
$ws.Range("I74").Value = 111237480
$ws.Range("K74").Value = 111237480
$ws.Range("M74").Value = -111236606

$ws.Range("H77").Value = 111237480
$ws.Range("I77").Value = 111237480
$ws.Range("K77").Value = 556187400
$ws.Range("M77").Value = -556183032

$ws.Range("H88").Value = 2448.3333
$ws.Range("J88").Value = 2466
$ws.Range("L88").Value = 2466
$ws.Range("N88").Value = -3278

$ws.Range("H91").Value = 2448.3333
$ws.Range("J91").Value = 2466
$ws.Range("L91").Value = 2466
$ws.Range("N91").Value = -5274

$ws.Range("H110").Value = 26707.934
$ws.Range("I110").Value = 28258.5
$ws.Range("K110").Value = 28258.5
$ws.Range("M110").Value = -26213.5

$ws.Range("H116").Value = 615
$ws.Range("I116").Value = 615
$ws.Range("K116").Value = 615
$ws.Range("M116").Value = 1679

$ws.Range("H122").Value = 3398.4167
$ws.Range("I122").Value = 2213.5
$ws.Range("J122").Value = 4583.3335
$ws.Range("K122").Value = 6640.5
$ws.Range("L122").Value = 13750.0005
$ws.Range("M122").Value = -4190.5
$ws.Range("N122").Value = -18650.0005

$ws.Range("H136").Value = 58828620
$ws.Range("I136").Value = 83336170
$ws.Range("K136").Value = 250008510
$ws.Range("M136").Value = -250005960

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 615
$ws.Range("I3").Value = 615
$ws.Range("K3").Value = 615
$ws.Range("M3").Value = -501

$ws.Range("H4").Value = 213.11111
$ws.Range("I4").Value = 136
$ws.Range("J4").Value = 367.33334
$ws.Range("K4").Value = 136
$ws.Range("L4").Value = 367.33334
$ws.Range("M4").Value = -21
$ws.Range("N4").Value = -597.33334

$ws.Range("H81").Value = 59424
$ws.Range("J81").Value = 59424
$ws.Range("L81").Value = 59424
$ws.Range("N81").Value = -61546

$ws.Range("H84").Value = 59424
$ws.Range("J84").Value = 59424
$ws.Range("L84").Value = 178272
$ws.Range("N84").Value = -188880

$ws.Range("H86").Value = 22371.334
$ws.Range("I86").Value = 5801.3335
$ws.Range("J86").Value = 27894.666
$ws.Range("K86").Value = 5801.3335
$ws.Range("L86").Value = 27894.666
$ws.Range("M86").Value = -4678.3335
$ws.Range("N86").Value = -30140.666

$ws.Range("H89").Value = 22371.334
$ws.Range("I89").Value = 5801.3335
$ws.Range("J89").Value = 27894.666
$ws.Range("K89").Value = 29006.6675
$ws.Range("L89").Value = 139473.33
$ws.Range("M89").Value = -23390.6675
$ws.Range("N89").Value = -150705.33

$ws.Range("H94").Value = 1265.75
$ws.Range("I94").Value = 1205.0385
$ws.Range("K94").Value = 1205.0385
$ws.Range("M94").Value = -754.0385000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 353.35
$ws.Range("I7").Value = 106.3
$ws.Range("J7").Value = 600.4
$ws.Range("K7").Value = 106.3
$ws.Range("L7").Value = 600.4
$ws.Range("M7").Value = 6.700000000000003
$ws.Range("N7").Value = -826.4

$ws.Range("H16").Value = 2709.9
$ws.Range("I16").Value = 3097.6667
$ws.Range("K16").Value = 3097.6667
$ws.Range("M16").Value = -2810.6667

$ws.Range("H22").Value = 3938.5518
$ws.Range("I22").Value = 5912.6665
$ws.Range("J22").Value = 708.1818
$ws.Range("K22").Value = 5912.6665
$ws.Range("L22").Value = 708.1818
$ws.Range("M22").Value = -5562.6665
$ws.Range("N22").Value = -1408.1818

$ws.Range("H31").Value = 34488316
$ws.Range("I31").Value = 3850.8
$ws.Range("K31").Value = 3850.8
$ws.Range("M31").Value = -3555.8

$ws.Range("H34").Value = 34488316
$ws.Range("I34").Value = 3850.8
$ws.Range("K34").Value = 3850.8
$ws.Range("M34").Value = -3648.8

$ws.Range("H58").Value = 1455.48
$ws.Range("I58").Value = 1402.8182
$ws.Range("K58").Value = 1402.8182
$ws.Range("M58").Value = -1199.8182

$ws.Range("H69").Value = 72567.71000000001
$ws.Range("I69").Value = 46594.8
$ws.Range("J69").Value = 137500
$ws.Range("K69").Value = 46594.8
$ws.Range("L69").Value = 137500
$ws.Range("M69").Value = -45845.8
$ws.Range("N69").Value = -138998

$ws.Range("H72").Value = 72567.71000000001
$ws.Range("I72").Value = 46594.8
$ws.Range("J72").Value = 137500
$ws.Range("K72").Value = 139784.4
$ws.Range("L72").Value = 412500
$ws.Range("M72").Value = -136040.4
$ws.Range("N72").Value = -419988

$ws.Range("H99").Value = 14519
$ws.Range("I99").Value = 17827.2
$ws.Range("K99").Value = 17827.2
$ws.Range("M99").Value = -16329.2

$ws.Range("H113").Value = 2709.9
$ws.Range("I113").Value = 3097.6667
$ws.Range("K113").Value = 3097.6667
$ws.Range("M113").Value = -927.6667000000002

$ws.Range("H126").Value = 14519
$ws.Range("I126").Value = 17827.2
$ws.Range("K126").Value = 53481.60000000001
$ws.Range("M126").Value = -51011.60000000001

$ws.Range("H132").Value = 3097.5151
$ws.Range("I132").Value = 2711.963
$ws.Range("J132").Value = 4832.5
$ws.Range("K132").Value = 8135.889000000001
$ws.Range("L132").Value = 14497.5
$ws.Range("M132").Value = -5605.889000000001
$ws.Range("N132").Value = -19557.5

$ws.Range("H134").Value = 1061.2307
$ws.Range("I134").Value = 974.8333
$ws.Range("J134").Value = 2098
$ws.Range("K134").Value = 2924.4999
$ws.Range("L134").Value = 6294
$ws.Range("M134").Value = -389.4998999999998
$ws.Range("N134").Value = -11364

$ws.Range("H136").Value = 1455.48
$ws.Range("I136").Value = 1402.8182
$ws.Range("K136").Value = 4208.4546
$ws.Range("M136").Value = -1658.4546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1665664.8
$ws.Range("J2").Value = 2422036.8
$ws.Range("L2").Value = 14532220.8
$ws.Range("N2").Value = -14532446.8

$ws.Range("H4").Value = 55300376
$ws.Range("I4").Value = 101437940
$ws.Range("J4").Value = 24542000
$ws.Range("K4").Value = 304313820
$ws.Range("L4").Value = 73626000
$ws.Range("M4").Value = -304313708
$ws.Range("N4").Value = -73626224

$ws.Range("H25").Value = 2385
$ws.Range("J25").Value = 5600
$ws.Range("L25").Value = 16800
$ws.Range("N25").Value = -17138

$ws.Range("H30").Value = 2385
$ws.Range("J30").Value = 5600
$ws.Range("L30").Value = 16800
$ws.Range("N30").Value = -17004

$ws.Range("H107").Value = 898.5789
$ws.Range("J107").Value = 954.8570999999999
$ws.Range("L107").Value = 2864.5713
$ws.Range("N107").Value = -6704.5713

$ws.Range("H131").Value = 1703.3
$ws.Range("J131").Value = 1783.2858
$ws.Range("L131").Value = 5349.857400000001
$ws.Range("N131").Value = -15429.8574

$ws.Range("H132").Value = 3336910.2
$ws.Range("I132").Value = 1916.6666
$ws.Range("J132").Value = 3925438.5
$ws.Range("K132").Value = 17249.9994
$ws.Range("L132").Value = 35328946.5
$ws.Range("M132").Value = -14719.9994
$ws.Range("N132").Value = -35334006.5

$ws.Range("H133").Value = 12985.889
$ws.Range("J133").Value = 19954.75
$ws.Range("L133").Value = 59864.25
$ws.Range("N133").Value = -69984.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1681.3636
$ws.Range("I97").Value = 1110.5555
$ws.Range("K97").Value = 1110.5555
$ws.Range("M97").Value = -614.5554999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2675.25
$ws.Range("I16").Value = 2675.25
$ws.Range("K16").Value = 2675.25
$ws.Range("M16").Value = -2505.25

$ws.Range("H22").Value = 3549.6365
$ws.Range("J22").Value = 4270.857
$ws.Range("L22").Value = 4270.857
$ws.Range("N22").Value = -4860.857

$ws.Range("H27").Value = 3549.6365
$ws.Range("J27").Value = 4270.857
$ws.Range("L27").Value = 4270.857
$ws.Range("N27").Value = -4484.857

$ws.Range("H46").Value = 1786.6842

$ws.Range("H82").Value = 3365.1667
$ws.Range("I82").Value = 1175.6666
$ws.Range("J82").Value = 5554.6665
$ws.Range("K82").Value = 1175.6666
$ws.Range("L82").Value = 5554.6665
$ws.Range("M82").Value = -814.6666
$ws.Range("N82").Value = -6276.6665

$ws.Range("H85").Value = 3365.1667
$ws.Range("I85").Value = 1175.6666
$ws.Range("J85").Value = 5554.6665
$ws.Range("K85").Value = 1175.6666
$ws.Range("L85").Value = 5554.6665
$ws.Range("M85").Value = 72.33339999999998
$ws.Range("N85").Value = -8050.6665

$ws.Range("H136").Value = 5140.839
$ws.Range("I136").Value = 4532.7407
$ws.Range("J136").Value = 9245.5
$ws.Range("K136").Value = 13598.2221
$ws.Range("L136").Value = 27736.5
$ws.Range("M136").Value = -11048.2221
$ws.Range("N136").Value = -32836.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 6072.636
$ws.Range("J96").Value = 7457
$ws.Range("L96").Value = 7457
$ws.Range("N96").Value = -10203

$ws.Range("H107").Value = 537.1667
$ws.Range("I107").Value = 584.25
$ws.Range("J107").Value = 443
$ws.Range("K107").Value = 1752.75
$ws.Range("L107").Value = 1329
$ws.Range("M107").Value = 167.25
$ws.Range("N107").Value = -5169

$ws.Range("H132").Value = 4825.2085
$ws.Range("I132").Value = 4557.8535
$ws.Range("J132").Value = 6391.143
$ws.Range("K132").Value = 13673.5605
$ws.Range("L132").Value = 19173.429
$ws.Range("M132").Value = -11143.5605
$ws.Range("N132").Value = -24233.429

$ws.Range("H137").Value = 100000
$ws.Range("J137").Value = 100000
$ws.Range("L137").Value = 100000
$ws.Range("N137").Value = -110200
